$data = @{
    2 = @(0.01514828764759746, 0.002777888934908601, 3.900430680208489, 0.496779210170732, 4.415136066961727)
    3 = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 3.645393585217082)
    4 = @(0.01514828764759746, 0.002777888934908601, 3.900430680208489, 0.496779210170732, 4.415136066961727)
    5 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    6 = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.997429241610044)
    7 = @(0.3048080303191223, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1.271902929317955)
    8 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 8.660232485948974, 17.45944343273191)
    9 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    10 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759)
    11 = @(0.3048080303191223, 0.3127903958511391, 3.900430680208489, 0.496779210170732, 5.014808316549482)
    12 = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.997429241610044)
    13 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    14 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797)
    15 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    16 = @(0.6753301551942219, 1.667794583268128, 3.900430680208489, 0.496779210170732, 6.740334628841572)
    17 = @(1.459612070389937, 1.667794583268128, 3.900430680208489, 0.496779210170732, 7.524616544037286)
    18 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    19 = @(0.6753301551942219, 1.667794583268128, 3.900430680208489, 0.496779210170732, 6.740334628841572)
    20 = @(0.127881588408715, 0.002777888934908601, 0.1575252929769615, 0.496779210170732, 0.7849639804913171)
    21 = @(0.04763786555579896, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1.014732764554632)
    22 = @(0.003994804209775715, 0.3127903958511391, 3.900430680208489, 0.496779210170732, 4.713995090440136)
    23 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    24 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    25 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    26 = @(1.459612070389937, 1.667794583268128, 3.900430680208489, 0.496779210170732, 7.524616544037286)
    27 = @(0.3048080303191223, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 1.919867272924993)
    28 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759)
    29 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    30 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    31 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759)
    32 = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 3.645393585217082)
    33 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759)
    34 = @(0.6753301551942219, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 2.290389397800092)
    35 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    36 = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.997429241610044)
    37 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    38 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    39 = @(0.127881588408715, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 1.742940831014585)
    40 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 4.429675500412797)
    41 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
